# Generate Report for Handback
# Update the timestamp cells on the Overview / zh-cn / de-de sheets to
# reflect the latest handback/handoff generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for 58ea362d-e041-4007-a48b-2f212891fbf1.md
$wsOverview.Range("G2").Value = "2016-08-17 06:58:29"

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for the same file
$wsZhCn.Range("H2").Value = "2016-08-17 06:58:22"
$wsZhCn.Range("K2").Value = "2016-08-17 06:58:41"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime for the same file
$wsDeDe.Range("H2").Value = "2016-08-17 06:58:29"
$wsDeDe.Range("K2").Value = "2016-08-17 06:58:48"
